# Applies the diff: adds "Project"/"Ilias" component sub-sections (Heading5 +
# descriptive paragraphs) at four locations in the SI (System and Information
# Integrity) control narrative. Insertions are performed from the bottom of
# the document upward so that earlier insertions don't shift the paragraph
# indices used by later ones.

$d = $word.ActiveDocument

# Positional params: this runtime's PowerShell subset does not bind
# `-ParamName value` style named arguments reliably, so use positional ones.
#
# NOTE: Paragraph.Next()/.Previous() do not return usable live references in
# this host (subsequent property reads/writes on them silently no-op), so
# every paragraph is instead (re)looked-up through $d.Paragraphs.Item(index)
# — that binding is reliable. AnchorIndex is the 1-based index of the
# paragraph the new one should be inserted AFTER; returns the new
# paragraph's 1-based index.
function Add-HeadingParagraph {
    param([int]$AnchorIndex, [string]$StyleName, [string]$Text, [string]$BookmarkName)
    $anchorPara = $d.Paragraphs.Item($AnchorIndex)
    $anchorPara.Range.InsertParagraphAfter() | Out-Null
    $newIndex = $AnchorIndex + 1
    $newPara = $d.Paragraphs.Item($newIndex)
    $newPara.Style = $StyleName
    $newPara.Range.Text = $Text
    if ($BookmarkName) {
        $newPara = $d.Paragraphs.Item($newIndex)
        $bmRange = $d.Range($newPara.Range.Start, $newPara.Range.Start + $Text.Length)
        $d.Bookmarks.Add($BookmarkName, $bmRange) | Out-Null
    }
    return $newIndex
}

# ---------------------------------------------------------------------------
# Insertion 4 (bottommost): end of document, after the final "CivicActions"
# SI-12 retention paragraph -> add "Project" Heading5 + FirstParagraph.
# ---------------------------------------------------------------------------
$n = $d.Paragraphs.Count

$idx = Add-HeadingParagraph $n "Heading 5" "Project" "project-2"
$idx = Add-HeadingParagraph $idx "First Paragraph" "Project representatives and systems administrators receive annual training from Client regarding information assurance and information handling requirements. These personnel are required to operate the system and handle system data and output in accordance with legal requirements. Personnel training and system guidelines ensure that data and programs are handled appropriately." $null

# ---------------------------------------------------------------------------
# Insertion 3: inside SI-5 section, right before the "a" Heading4 sub-item
# (after the SourceCode + "Status: Complete" paragraphs) -> add "Ilias" and
# "Project" Heading5 blocks.
# ---------------------------------------------------------------------------
$n = $d.Paragraphs.Count
$targetIdx = 0
$foundSI5 = $false
for ($i = 1; $i -le $n; $i++) {
    $para = $d.Paragraphs.Item($i)
    if ($para.Range.Text -like "*SI-5: Security Alerts*") {
        $foundSI5 = $true
    }
    if ($foundSI5 -and $para.Style.NameLocal -eq "Heading 4" -and $para.Range.Text.Trim() -eq "a") {
        $targetIdx = $i
        break
    }
}
$anchorIdx = $targetIdx - 1  # "Status: Complete" paragraph

$idx = Add-HeadingParagraph $anchorIdx "Heading 5" "Ilias" "ilias-1"
$idx = Add-HeadingParagraph $idx "First Paragraph" "CivicActions Security and Operations receive Ilias Security Advisories on a regular basis." $null
$idx = Add-HeadingParagraph $idx "Heading 5" "Project" "project-1"
$idx = Add-HeadingParagraph $idx "First Paragraph" "Project representatives and system administrators receive alerts from US-CERT on a regular basis. Support personnel take appropriate action in response to relevant areas of concern." $null

# ---------------------------------------------------------------------------
# Insertion 2: inside SI-2 section, right before the "a" Heading4 sub-item
# (after the SourceCode + "Status: Complete" paragraphs) -> add "Ilias"
# Heading5 block.
# ---------------------------------------------------------------------------
$n = $d.Paragraphs.Count
$targetIdx = 0
$foundSI2 = $false
for ($i = 1; $i -le $n; $i++) {
    $para = $d.Paragraphs.Item($i)
    if ($para.Range.Text -like "*SI-2: Flaw Remediation*") {
        $foundSI2 = $true
    }
    if ($foundSI2 -and $para.Style.NameLocal -eq "Heading 4" -and $para.Range.Text.Trim() -eq "a") {
        $targetIdx = $i
        break
    }
}
$anchorIdx = $targetIdx - 1  # "Status: Complete" paragraph

$idx = Add-HeadingParagraph $anchorIdx "Heading 5" "Ilias" "ilias"
$idx = Add-HeadingParagraph $idx "First Paragraph" "Ilias contains built-in security status monitoring of the core application and contributed modules." $null

# ---------------------------------------------------------------------------
# Insertion 1 (topmost): right before the "SI-2: Flaw Remediation" Heading3
# -> add "Project" Heading5 + FirstParagraph + two BodyText paragraphs.
# ---------------------------------------------------------------------------
$n = $d.Paragraphs.Count
$targetIdx = 0
for ($i = 1; $i -le $n; $i++) {
    $para = $d.Paragraphs.Item($i)
    if ($para.Style.NameLocal -eq "Heading 3" -and $para.Range.Text.Trim() -eq "SI-2: Flaw Remediation") {
        $targetIdx = $i
        break
    }
}
$anchorIdx = $targetIdx - 1  # CivicActions policy FirstParagraph

$idx = Add-HeadingParagraph $anchorIdx "Heading 5" "Project" "project"
$idx = Add-HeadingParagraph $idx "First Paragraph" "System and information integrity policy and procedures for the Project system are formally documented in the Project SSP, which provides the roles and responsibilities as it pertains to physical and environmental protection systems. The Project system support staff monitors the network on a daily basis and employs up-to-date patches to protect the integrity of the system." $null
$idx = Add-HeadingParagraph $idx "Body Text" "Additional information is contained within the None." $null
$idx = Add-HeadingParagraph $idx "Body Text" "This is Agency common control. More data about implementation can be obtained from the Agency common control catalog." $null

Write-Output "final paragraph count: $($d.Paragraphs.Count)"
